# Further work on BDO and lipidcane biodiesel pretreatment:
# - Rename the "Additional lipid extraction efficiency [%]" parameter label
#   to "Bagasse lipid extraction efficiency [%]".
# - Refresh the Spearman correlation coefficients / p-values table (rows 4-13,
#   columns C, D, E, F, H, I, J) with updated results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the parameter label in B6.
$ws.Range("B6").Value = "Bagasse lipid extraction efficiency [%]"

# 2. Update the numeric results table.
$values = @{
    4  = @{ C = -0.02039300982524563;  D = 0.9029070726768171;   E = -0.9732258306457664;  F = 0.8416905422635568;  H = 0.1085607140178505;   I = 0.6205175129378234;   J = 0.7601530038250958 }
    5  = @{ C = 0.1104372609315233;    D = 0.2149898747468687;   E = -0.04046801170029252; F = 0.1897907447686192;  H = 0.2767329183229582;   I = 0.1823130578264457;   J = 0.1504582614565364 }
    6  = @{ C = 0.1106262656566414;    D = 0.337094427360684;    E = -0.01311032775819396; F = -0.1315667891697292; H = 0.05746793669841747;  I = -0.1340403510087753;  J = 0.1284437110927774 }
    7  = @{ C = 0.1360939023475587;    D = 0.08491862296557415;  E = 0.2463901597539939;   F = 0.2850956273906848;  H = 0.7307012675316885;   I = 0.6927743193579841;   J = -0.1410440261006525 }
    8  = @{ C = 0.6631080777019427;    D = 0.0461921548038701;   E = -0.02728568214205356; F = 0.06711617790444763; H = 0.04327608190204756;  I = 0.06610965274131855;  J = 0.03392934823370585 }
    9  = @{ C = 0.4182374559363984;    D = -0.09282082052051303; E = 0.07926798169954251;  F = -0.0819545488637216; H = -0.02108602715067877; I = -0.02833570839270982; J = -0.0726558163954099 }
    10 = @{ C = 0.102977574439361;     D = -0.002338558463961599; E = -0.04251406285157129; F = 0.02496662416560414; H = 0.01120978024450611;  I = -0.006360159003975101; J = 0.03134028350708768 }
    11 = @{ C = 0.192606315157879;     D = 0.063841596039901;    E = -0.03549688742218556; F = 0.07115577889447237; H = 0.08450311257781445;  I = 0.1063811595289882;   J = 0.00633015825395635 }
    12 = @{ C = 0.1868896722418061;    D = 0.1116717917947949;   E = 0.1335498387459687;   F = 0.2536038400960025;  H = 0.5272991824795621;   I = -0.06806570164254108; J = 0.003169579239480987 }
    13 = @{ C = -0.3846066151653792;   D = -0.001168529213230331; E = -0.008241206030150754; F = 0.01636690917272932; H = 0.0086987174679367;   I = 0.05195979899497488;  J = -0.0214310357758944 }
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}
